# Emilia: Kleine Änderung Powerpoint
#
# Swaps the title text and the picture geometry between slide 13
# ("Feigenbaum-Diagramm: Warum dieser Bereich von r?") and slide 14
# ("Logistische Rekursion mit 3 Häufungspunkten") so that each slide now
# shows the other slide's title/picture layout.

$p = $ppt.ActivePresentation

# Shape.Left/Top/Width/Height round-trip through a lower-precision
# float representation, so naively setting "targetEmu / 12700" as points
# can land the saved OOXML value a handful of EMU away from the exact
# target. Do a small local search (in EMU-sized steps, not point-sized
# steps) around the naive points value so we land on the exact EMU.
function Set-ExactLeft($shape, $targetEmu) {
    for ($i = -2000; $i -le 2000; $i++) {
        $pts = ($targetEmu + ($i * 0.01)) / 12700.0
        $shape.Left = $pts
        $gotEmu = [Math]::Round($shape.Left * 12700)
        if ($gotEmu -eq $targetEmu) {
            return
        }
    }
}

function Set-ExactTop($shape, $targetEmu) {
    for ($i = -2000; $i -le 2000; $i++) {
        $pts = ($targetEmu + ($i * 0.01)) / 12700.0
        $shape.Top = $pts
        $gotEmu = [Math]::Round($shape.Top * 12700)
        if ($gotEmu -eq $targetEmu) {
            return
        }
    }
}

function Set-ExactWidth($shape, $targetEmu) {
    for ($i = -2000; $i -le 2000; $i++) {
        $pts = ($targetEmu + ($i * 0.01)) / 12700.0
        $shape.Width = $pts
        $gotEmu = [Math]::Round($shape.Width * 12700)
        if ($gotEmu -eq $targetEmu) {
            return
        }
    }
}

function Set-ExactHeight($shape, $targetEmu) {
    for ($i = -2000; $i -le 2000; $i++) {
        $pts = ($targetEmu + ($i * 0.01)) / 12700.0
        $shape.Height = $pts
        $gotEmu = [Math]::Round($shape.Height * 12700)
        if ($gotEmu -eq $targetEmu) {
            return
        }
    }
}

$slide13 = $p.Slides.Item(13)
$slide14 = $p.Slides.Item(14)

$title13 = $slide13.Shapes.Item(1)
$pic13 = $slide13.Shapes.Item(2)

$title14 = $slide14.Shapes.Item(1)
$pic14 = $slide14.Shapes.Item(2)

# Swap the titles.
$title13.TextFrame.TextRange.Text = "Logistische Rekursion mit 3 Häufungspunkten"
$title14.TextFrame.TextRange.Text = "Feigenbaum-Diagramm: Warum dieser Bereich von r?"

# Swap the picture geometry (EMUs).
Set-ExactLeft   $pic13 1507163
Set-ExactTop    $pic13 1017725
Set-ExactWidth  $pic13 6129668
Set-ExactHeight $pic13 3820975

Set-ExactLeft   $pic14 1653200
Set-ExactTop    $pic14 1017725
Set-ExactWidth  $pic14 5837601
Set-ExactHeight $pic14 3820975

Write-Output "done"
